$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 13:50"

$data = @(
    @("Madrid", 24090, 8301, 12397, 3392),
    @("Cataluña", 15026, 3455, 10345, 1226),
    @("Bizkaia/Vizcaya", 2937, 1626, 2054, 133),
    @("Castilla-La Mancha", 2780, 71, 2446, 263),
    @("Valencia/Valencia", 2685, 130, 2412, 143),
    @("Navarra", 2146, 161, 1883, 102),
    @("Araba/Alava", 2049, 1626, 1376, 123),
    @("Ciudad Real", 2041, 252, 1772, 179),
    @("Alacant/Alicante", 1839, 46, 1658, 135),
    @("La Rioja", 1733, 431, 1231, 71),
    @("Zaragoza", 1641, 141, 1419, 81),
    @("Albacete", 1537, 252, 1339, 133),
    @("Toledo", 1426, 252, 1207, 155),
    @("A Coruña", 1351, 153, 1261, 40),
    @("Malaga", 1321, 83, 1169, 69),
    @("Asturias", 1158, 78, 1031, 49),
    @("Cantabria", 1100, 24, 1049, 27),
    @("Salamanca", 1078, 181, 787, 110),
    @("Gipuzkoa/Guipuzcoa", 1071, 1626, 684, 41),
    @("Pontevedra", 1060, 153, 1005, 9),
    @("Sevilla", 1052, 18, 1000, 34),
    @("Caceres", 1045, 11, 945, 89),
    @("Granada", 963, 15, 882, 66),
    @("Murcia", 939, 17, 897, 25),
    @("Valladolid", 929, 145, 718, 66),
    @("Leon", 918, 139, 690, 89),
    @("Aragon", 907, 29, 838, 40),
    @("Burgos", 749, 176, 512, 61),
    @("La Palma", 712, 30, 1056, 2),
    @("Segovia", 629, 156, 404, 69),
    @("Jaen", 599, 17, 559, 23),
    @("Castello/Castellon", 586, 9, 545, 32),
    @("Guadalajara", 586, 252, 479, 93),
    @("Cordoba", 572, 4, 555, 13),
    @("Soria", 550, 71, 442, 37),
    @("Badajoz", 515, 49, 449, 17),
    @("Cadiz", 507, 10, 484, 13),
    @("Ourense", 458, 153, 415, 8),
    @("Avila", 446, 91, 309, 46),
    @("Palencia", 293, 33, 242, 18),
    @("Fuerteventura", 288, 30, 1056, 0),
    @("Lugo", 270, 153, 244, 4),
    @("Cuenca", 268, 252, 187, 62),
    @("Almeria", 223, 6, 203, 14),
    @("Teruel", 222, 14, 196, 12),
    @("Huesca", 215, 19, 185, 11),
    @("Mallorca", 210, 18, 194, 12),
    @("Zamora", 209, 36, 153, 20),
    @("Huelva", 168, 2, 162, 4),
    @("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @("Melilla", 51, 0, 50, 1),
    @("Lanzarote", 49, 30, 1056, 3),
    @("La Gomera", 35, 30, 1056, 0),
    @("El Hierro", 30, 30, 1056, 0),
    @("Ceuta", 29, 0, 28, 1),
    @("Ibiza", 21, 18, 20, 1),
    @("Menorca", 15, 18, 13, 0),
    @("Gran Canaria", 8, 30, 1056, 11),
    @("Arroyo de la Luz", 7, 0, 7, 0),
    @("Tenerife", 3, 30, 1056, 36),
    @("Formentera", 0, 10, 0, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
